# Generate Report for Handoff
# Adds two new handed-off files (6eb1c6f5-... and c34fd46a-...) as new rows
# to the "Overview", "zh-cn" and "de-de" sheets, extending their tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New file identifiers
# ---------------------------------------------------------------------------
$guid1 = "6eb1c6f5-0f9d-4b6c-8e7f-161097ffd532"
$guid2 = "c34fd46a-8e3b-4ad9-9a3e-87d6ed99e835"

$md1 = "$guid1.md"
$md2 = "$guid2.md"

$path1 = "e2e\$md1"
$path2 = "e2e\$md2"

$dateOverview = "2016-08-14 02:58:03"
$dateZh = "2016-08-14 02:57:53"
$dateDe = "2016-08-14 02:58:03"

$xlfZh1 = "$guid1.400a744314e8f69fc5b3b6e8a4ff8ad664f46aea.zh-cn.xlf"
$xlfZh2 = "$guid2.65686b608a8c30f7c8108353c7187a67a7f93d23.zh-cn.xlf"
$xlfDe1 = "$guid1.400a744314e8f69fc5b3b6e8a4ff8ad664f46aea.de-de.xlf"
$xlfDe2 = "$guid2.65686b608a8c30f7c8108353c7187a67a7f93d23.de-de.xlf"

$urlBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob"

# ---------------------------------------------------------------------------
# Overview sheet: two new rows (6 & 7), columns A-G
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A6").Value = $md1
$wsOverview.Range("B6").Value = $path1
$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("D6").Value = ""
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = $dateOverview
$wsOverview.Range("G6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Range("A7").Value = $md2
$wsOverview.Range("B7").Value = $path2
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("D7").Value = ""
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = $dateOverview
$wsOverview.Range("G7").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "$urlBase/0000000000000000000000000000000000000000/e2e/$md1", "", "", $path1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "$urlBase/0000000000000000000000000000000000000000/e2e/$md2", "", "", $path2)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G7"))

# ---------------------------------------------------------------------------
# zh-cn sheet: two new rows (6 & 7), columns A-P
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A6").Value = $md1
$wsZh.Range("B6").Value = ".md"
$wsZh.Range("C6").Value = "Ready for handoff"
$wsZh.Range("D6").Value = "e2e"
$wsZh.Range("E6").Value = "ht"
$wsZh.Range("F6").Value = "'False"
$wsZh.Range("G6").Value = $xlfZh1
$wsZh.Range("H6").Value = $dateZh
$wsZh.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I6").Value = ""
$wsZh.Range("J6").Value = ""
$wsZh.Range("K6").Value = "0001-01-01 00:00:00"
$wsZh.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L6").Value = ""
$wsZh.Range("M6").Value = "'True"
$wsZh.Range("N6").Value = ""
$wsZh.Range("O6").Value = "'False"
$wsZh.Range("P6").Value = ""

$wsZh.Range("A7").Value = $md2
$wsZh.Range("B7").Value = ".md"
$wsZh.Range("C7").Value = "Ready for handoff"
$wsZh.Range("D7").Value = "e2e"
$wsZh.Range("E7").Value = "ht"
$wsZh.Range("F7").Value = "'False"
$wsZh.Range("G7").Value = $xlfZh2
$wsZh.Range("H7").Value = $dateZh
$wsZh.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I7").Value = ""
$wsZh.Range("J7").Value = ""
$wsZh.Range("K7").Value = "0001-01-01 00:00:00"
$wsZh.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L7").Value = ""
$wsZh.Range("M7").Value = "'True"
$wsZh.Range("N7").Value = ""
$wsZh.Range("O7").Value = "'False"
$wsZh.Range("P7").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "$urlBase/0000000000000000000000000000000000000000/e2e/$md1", "", "", $md1)
$wsZh.Hyperlinks.Add($wsZh.Range("A7"), "$urlBase/0000000000000000000000000000000000000000/e2e/$md2", "", "", $md2)

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P7"))

# ---------------------------------------------------------------------------
# de-de sheet: two new rows (6 & 7), columns A-P
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A6").Value = $md1
$wsDe.Range("B6").Value = ".md"
$wsDe.Range("C6").Value = "Ready for handoff"
$wsDe.Range("D6").Value = "e2e"
$wsDe.Range("E6").Value = "ht"
$wsDe.Range("F6").Value = "False"
$wsDe.Range("G6").Value = $xlfDe1
$wsDe.Range("H6").Value = $dateDe
$wsDe.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I6").Value = ""
$wsDe.Range("J6").Value = ""
$wsDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDe.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L6").Value = ""
$wsDe.Range("M6").Value = "True"
$wsDe.Range("N6").Value = ""
$wsDe.Range("O6").Value = "False"
$wsDe.Range("P6").Value = ""

$wsDe.Range("A7").Value = $md2
$wsDe.Range("B7").Value = ".md"
$wsDe.Range("C7").Value = "Ready for handoff"
$wsDe.Range("D7").Value = "e2e"
$wsDe.Range("E7").Value = "ht"
$wsDe.Range("F7").Value = "False"
$wsDe.Range("G7").Value = $xlfDe2
$wsDe.Range("H7").Value = $dateDe
$wsDe.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I7").Value = ""
$wsDe.Range("J7").Value = ""
$wsDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDe.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L7").Value = ""
$wsDe.Range("M7").Value = "True"
$wsDe.Range("N7").Value = ""
$wsDe.Range("O7").Value = "False"
$wsDe.Range("P7").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "$urlBase/0000000000000000000000000000000000000000/e2e/$md1", "", "", $md1)
$wsDe.Hyperlinks.Add($wsDe.Range("A7"), "$urlBase/0000000000000000000000000000000000000000/e2e/$md2", "", "", $md2)

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P7"))
